$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4050.423
$ws.Range("C2").Value = 3795.6358
$ws.Range("D2").Value = 372.9026
$ws.Range("E2").Value = 11936.176
$ws.Range("F2").Value = 11806.5536
$ws.Range("G2").Value = 133.3958

$ws.Range("B3").Value = 5024.7916
$ws.Range("C3").Value = 4749.8566
$ws.Range("D3").Value = 369.4052
$ws.Range("E3").Value = 11886.532
$ws.Range("F3").Value = 11783.1164
$ws.Range("G3").Value = 124.9352

$ws.Range("B4").Value = 4328.3129
$ws.Range("C4").Value = 3796.8726
$ws.Range("D4").Value = 636.2219
$ws.Range("E4").Value = 15525.2284
$ws.Range("F4").Value = 14363.384
$ws.Range("G4").Value = 1238.3755

$ws.Range("B5").Value = 4215.04
$ws.Range("C5").Value = 3801.8807
$ws.Range("D5").Value = 476.3175
$ws.Range("E5").Value = 29336.968
$ws.Range("F5").Value = 29284.2368
$ws.Range("G5").Value = 152.1815

$ws.Range("B6").Value = 5191.954
$ws.Range("C6").Value = 4743.1756
$ws.Range("D6").Value = 472.0929
$ws.Range("E6").Value = 29335.7544
$ws.Range("F6").Value = 29347.1344
$ws.Range("G6").Value = 133.4588

$ws.Range("B7").Value = 4348.4539
$ws.Range("C7").Value = 3759.0585
$ws.Range("D7").Value = 678.572
$ws.Range("E7").Value = 35878.1336
$ws.Range("F7").Value = 35176.5512
$ws.Range("G7").Value = 784.7686

$ws.Range("B8").Value = 4160.6541
$ws.Range("C8").Value = 3695.0979
$ws.Range("D8").Value = 502.1258
$ws.Range("E8").Value = 29284.9672
$ws.Range("F8").Value = 29273.4584
$ws.Range("G8").Value = 167.6822

$ws.Range("B9").Value = 4241.7386
$ws.Range("C9").Value = 3806.5155
$ws.Range("D9").Value = 491.6633
$ws.Range("E9").Value = 29372.6688
$ws.Range("F9").Value = 29242.7488
$ws.Range("G9").Value = 144.8376

$ws.Range("B10").Value = 4982.7126
$ws.Range("C10").Value = 3867.8875
$ws.Range("D10").Value = 1184.4444
$ws.Range("E10").Value = 36150.96
$ws.Range("F10").Value = 34313.0208
$ws.Range("G10").Value = 1843.5925

$ws.Range("B11").Value = 5437.8142
$ws.Range("C11").Value = 4745.7178
$ws.Range("D11").Value = 778.2706
$ws.Range("E11").Value = 29373.2224
$ws.Range("F11").Value = 29209.6776
$ws.Range("G11").Value = 221.4471

$ws.Range("B12").Value = 5242.876
$ws.Range("C12").Value = 4757.2494
$ws.Range("D12").Value = 597.5981
$ws.Range("E12").Value = 29361.7464
$ws.Range("F12").Value = 29155.2184
$ws.Range("G12").Value = 167.3455

$ws.Range("B13").Value = 4335.62
$ws.Range("C13").Value = 3709.8307
$ws.Range("D13").Value = 698.2141
$ws.Range("E13").Value = 29401.6152
$ws.Range("F13").Value = 29133.548
$ws.Range("G13").Value = 184.6295

$ws.Range("B14").Value = 5288.4458
$ws.Range("C14").Value = 4759.7288
$ws.Range("D14").Value = 622.655
$ws.Range("E14").Value = 67147.0112
$ws.Range("F14").Value = 66817.416
$ws.Range("G14").Value = 343.3874

$ws.Range("B15").Value = 5428.132
$ws.Range("C15").Value = 4761.335
$ws.Range("D15").Value = 769.1126
$ws.Range("E15").Value = 67757.6912
$ws.Range("F15").Value = 67653.2304
$ws.Range("G15").Value = 246.9598
